# Applies price (D) and volume% (E) updates to the crypto symbol list
# per the Jan 9 2023 GitHub Actions data refresh.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'279.44"
$ws.Range("E2").Value = "'6.86%"
$ws.Range("D3").Value = "'27.40"
$ws.Range("E3").Value = "'2.48%"
$ws.Range("D4").Value = "'4.802"
$ws.Range("E4").Value = "'2.14%"
$ws.Range("D5").Value = "'0.06341"
$ws.Range("E5").Value = "'2.38%"
$ws.Range("D6").Value = "'6.944"
$ws.Range("E6").Value = "'2.90%"
$ws.Range("D7").Value = "'3.386"
$ws.Range("E7").Value = "'6.70%"
$ws.Range("D8").Value = "'0.8799"
$ws.Range("E8").Value = "'3.64%"
$ws.Range("D9").Value = "'0.9420"
$ws.Range("E9").Value = "'3.31%"
$ws.Range("D10").Value = "'0.1471"
$ws.Range("E10").Value = "'4.83%"
$ws.Range("D11").Value = "'0.05150"
$ws.Range("E11").Value = "'1.04%"
$ws.Range("D12").Value = "'0.07325"
$ws.Range("E12").Value = "'3.17%"
$ws.Range("D13").Value = "'0.03137"
$ws.Range("E13").Value = "'0.97%"
$ws.Range("D14").Value = "'0.09073"
$ws.Range("E14").Value = "'0.35%"
$ws.Range("D15").Value = "'0.001558"
$ws.Range("E15").Value = "'1.73%"
$ws.Range("D16").Value = "'0.0006271"
$ws.Range("E16").Value = "'1.75%"
$ws.Range("D17").Value = "'0.006050"
$ws.Range("E17").Value = "'1.78%"
$ws.Range("D18").Value = "'3.443"
$ws.Range("E18").Value = "'-0.11%"
$ws.Range("E19").Value = "'4.76%"
$ws.Range("E20").Value = "'2.64%"
$ws.Range("E21").Value = "'0.09%"
$ws.Range("D22").Value = "'3.860"
$ws.Range("E22").Value = "'-5.92%"
$ws.Range("D23").Value = "'0.04336"
$ws.Range("E23").Value = "'1.88%"
$ws.Range("D24").Value = "'0.001180"
$ws.Range("E24").Value = "'-0.34%"
$ws.Range("D25").Value = "'0.004303"
$ws.Range("E25").Value = "'6.02%"
$ws.Range("D26").Value = "'0.0001199"
$ws.Range("E26").Value = "'-0.17%"
$ws.Range("D27").Value = "'0.0001688"
$ws.Range("E27").Value = "'2.97%"
$ws.Range("D40").Value = "'0.04086"
$ws.Range("E40").Value = "'3.01%"
$ws.Range("D41").Value = "'0.006680"
$ws.Range("E41").Value = "'61.37%"
$ws.Range("D42").Value = "'0.1166"
$ws.Range("E42").Value = "'4.98%"
$ws.Range("D43").Value = "'0.002198"
$ws.Range("E43").Value = "'4.58%"
$ws.Range("D44").Value = "'0.01310"
$ws.Range("E44").Value = "'-1.39%"
$ws.Range("D45").Value = "'0.00005228"
$ws.Range("E45").Value = "'1.24%"
$ws.Range("E46").Value = "'-0.18%"
$ws.Range("E47").Value = "'854.37%"
$ws.Range("D48").Value = "'0.02248"
$ws.Range("E48").Value = "'-33.93%"
$ws.Range("D49").Value = "'0.00002098"
$ws.Range("E49").Value = "'-0.18%"
$ws.Range("D50").Value = "'0.0001998"
$ws.Range("E50").Value = "'-0.18%"
